$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 9.5
$ws.Range("B4").Value = 0.9
$ws.Range("C4").Value = 1.25
$ws.Range("C5").Value = 15

$ws.Range("C5").Select()
